# ADC System Overview.xlsx touch-ups:
#  - Rename "AXi" machine/outputs terminology to "AVI" (and ADC Outputs -> ADCS Outputs)
#  - Update the saved selection / zoom level left behind by the editing session

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Terminology touch-ups (shared-string text updates) ---
$ws.Range("D3").Value = "AVI Machine"
$ws.Range("F4").Value = "* AVI Outputs"
$ws.Range("J4").Value = "* ADCS Outputs"

# --- Window / view state ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("J5").Select()
